# StudyDesign.xlsx update: add "readme" sheet and "DosePerSurfaceArea" sheet.
$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheets in the order that reproduces the target sheetId
#        allocation (DosePerSurfaceArea first -> sheetId 2, readme second -> sheetId 3) ---
$doseBW = $wb.Worksheets.Item("DosePerBodyweight")

$newBSA = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $doseBW)
$newBSA.Name = "DosePerSurfaceArea"

$newReadme = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$newReadme.Name = "readme"

# Sheet handles above can become stale once the sheet collection is
# reshuffled (handles here track by slot/index, not identity) -- re-resolve
# every worksheet we still need to touch by name before using it further.
$doseBSA = $wb.Worksheets.Item("DosePerSurfaceArea")
$readme = $wb.Worksheets.Item("readme")

# --- 2. Populate "DosePerSurfaceArea" (mirrors DosePerBodyweight, BSA flavour) ---
$doseBSA.Range("A1").Value = "functionHandle = @addDosetablePerBSA"
$doseBSA.Range("A2").Value = "targetParameterList = {'*Application_*|ProtocolSchemaItem|DrugMass'}"
$doseBSA.Range("A3").Value = "BSAmin"
$doseBSA.Range("B3").Value = "BSAmax"
$doseBSA.Range("C3").Value = "targetParameter"
$doseBSA.Range("D3").Value = "dose_mg"

$doseBSA.Range("A4").Value = 50
$doseBSA.Range("B4").Value = 60
$doseBSA.Range("D4").Value = 55
$doseBSA.Range("C4").Formula = "=D4/225.21*1000"

$doseBSA.Range("A5").Formula = "=A4+10"
$doseBSA.Range("B5").Formula = "=B4+10"
$doseBSA.Range("D5").Formula = "=D4+10"
$doseBSA.Range("C5").Formula = "=D5/225.21*1000"

$doseBSA.Range("A6").Formula = "=A5+10"
$doseBSA.Range("B6").Formula = "=B5+10"
$doseBSA.Range("D6").Formula = "=D5+10"
$doseBSA.Range("C6").Formula = "=D6/225.21*1000"

$doseBSA.Range("A7").Formula = "=A6+10"
$doseBSA.Range("B7").Formula = "=B6+10"
$doseBSA.Range("D7").Formula = "=D6+10"
$doseBSA.Range("C7").Formula = "=D7/225.21*1000"

$doseBSA.Range("A8").Formula = "=A7+10"
$doseBSA.Range("B8").Formula = "=B7+10"
$doseBSA.Range("D8").Formula = "=D7+10"
$doseBSA.Range("C8").Formula = "=D8/225.21*1000"

$doseBSA.Range("A9").Formula = "=A8+10"
$doseBSA.Range("B9").Formula = "=B8+10"
$doseBSA.Range("D9").Formula = "=D8+10"
$doseBSA.Range("C9").Formula = "=D9/225.21*1000"

$doseBSA.Range("A10").Formula = "=A9+10"
$doseBSA.Range("B10").Formula = "=B9+10"
$doseBSA.Range("D10").Formula = "=D9+10"
$doseBSA.Range("C10").Formula = "=D10/225.21*1000"

# ColumnWidth is quantized to 1/7-character (pixel) steps by the host, same
# as real Excel -- use the pre-solved inputs that land exactly/closest on
# the template's stored widths (14 and 17.59765625 respectively).
$doseBSA.Columns.Item(1).ColumnWidth = 13.214285714285714
$doseBSA.Columns.Item(3).ColumnWidth = 16.785714285714285

$doseBSA.PageSetup.TopMargin = 56.69291339999997
$doseBSA.PageSetup.BottomMargin = 56.69291339999997

$doseBSA.Range("F11").Select()

# --- 3. Populate "readme" ---
$readme.Range("A1").Value = "First column function handel wich is used to set an application parameter
available are @addDosetablePerWeight and @addDosetablePerBSA"
$readme.Range("A2").Value = "second hadle list of parameters which are set by the function"
$readme.Range("A3").Value = "Third line header for numeric info, 
for addDosetablePerWeight folllwoing columns are mandatory column BWmin, BWmax and targetParameter
for @addDosetablePerBSA folllwoing columns are mandatory column BSAmin, BSAmax and targetParameter
"
$readme.Range("A4").Value = "additional colmuns can be added. Please make sure column name should no contain specialletters, also no spaces"
$readme.Range("A6").Value = "attention MoBi internla Units are used , Body weight [kk], BSA [dm^2], DrugMass µmol "

$readme.Range("A1").WrapText = $true
$readme.Range("A3").WrapText = $true

$readme.Rows.Item(1).RowHeight = 27.6
$readme.Rows.Item(3).RowHeight = 82.8

# Closest reachable width to the template's 72.19921875 (see ColumnWidth note above).
$readme.Columns.Item(1).ColumnWidth = 71.35714285714286

$readme.PageSetup.TopMargin = 56.69291339999997
$readme.PageSetup.BottomMargin = 56.69291339999997

$readme.Activate()
$readme.Range("A15").Select()
